$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("final_figures/2008/F1_P7_Hollensbe, Khazanchi & Masterson_2008_AMJ_How do I assess if my Supervisor and Organization are fair.png", "Yes", "Yes", "No", "No", "No", "Yes", "Yes"),
    @("final_figures/2008/F2_P19_Jarzabkowski_2008_AMJ_Shaping Strategy as a Structuration Process.png", "No", "Not applicable", "No", "No", "No", "No", "No"),
    @("final_figures/2008/F1_P27_O'Mahony & Bechky_2008_ASQ_Boundary Organizations enabling Collaboration among Unexpected Allies.png", "Yes", "Yes", "No", "No", "No", "Yes", "No"),
    @("final_figures/2008/F1_P10_Weber et al._2008_ASQ_Forage for Thought Mobilizing Codes in the movement for Grass-Fed Meat and Dairy Products.png", "Yes", "Yes", "No", "No", "No", "No", "No"),
    @("final_figures/2008/F2_P11_Weber et al._2008_ASQ_Forage for Thought Mobilizing Codes in the movement for Grass-Fed Meat and Dairy Products.png", "Yes", "Yes", "No", "No", "No", "No", "No"),
    @("final_figures/2008/F3_P12_Weber et al._2008_ASQ_Forage for Thought Mobilizing Codes in the movement for Grass-Fed Meat and Dairy Products.png", "Yes", "Yes", "No", "No", "No", "No", "No"),
    @("final_figures/2008/F4_P28_Weber et al._2008_ASQ_Forage for Thought Mobilizing Codes in the movement for Grass-Fed Meat and Dairy Products.png", "Yes", "Yes", "No", "No", "Yes", "Yes", "No"),
    @("final_figures/2008/F1_P3_Jarzabkowski_2008_AMJ_Shaping Strategy as a Structuration Process.png", "Yes", "Yes", "Yes", "No", "No", "No", "No")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
